$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 202, pushing existing rows 202-219 down to 203-220
$ws.Rows(202).Insert()

# Populate the newly inserted row 202 with the new weekly price record
$ws.Range("A202").Value = 7
$ws.Range("B202").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C202").Value = "Ñuble"
$ws.Range("D202").Value = 44769
$ws.Range("E202").Value = 16
$ws.Range("F202").Value = 100112017
$ws.Range("G202").Value = "Apio"
$ws.Range("H202").Value = "Americana (o)"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 80
$ws.Range("K202").Value = 9500
$ws.Range("L202").Value = 10000
$ws.Range("M202").Value = 9750
$ws.Range("N202").Value = "`$/docena de matas"
$ws.Range("O202").Value = "Provincia del Elquí"
$ws.Range("P202").Value = 1625
$ws.Range("Q202").Value = 6
$ws.Range("R202").Value = "Hortaliza"
